$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (shifts E onward to the right by one,
# and Excel copies the formatting of column D into the new column E).
$ws.Columns("E:E").Insert()

# Header row: D1 text changes from "# Deployments" to "Oct. Deployments";
# the new E1 cell becomes "Sep. Deployments".
$ws.Range("D1").Value = "Oct. Deployments"
$ws.Range("E1").Value = "Sep. Deployments"

# Fill in the new "Sep. Deployments" data column (rows 2-5).
$ws.Range("E2").Value = 4
$ws.Range("E3").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 6

# Match column E's width to column D's width (Excel normally copies the
# left-neighbor's width when inserting, rounded to character units).
$ws.Columns("E:E").ColumnWidth = $ws.Columns("D:D").ColumnWidth

# Move the active selection, as recorded in the saved workbook.
$ws.Range("E11").Select()
